$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.061.88'
$ws.Range("E2").Value = '  +3.80%  '
$ws.Range("D3").Value = '1.894.53'
$ws.Range("E3").Value = '  +4.11%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9980'
$ws.Range("E4").Value = '  -0.31%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '247.55'
$ws.Range("E5").Value = '  +0.56%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9982'
$ws.Range("E6").Value = '  -0.32%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4976'
$ws.Range("E7").Value = '  +1.16%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '44.78'
$ws.Range("E8").Value = '  +0.96%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2956'
$ws.Range("E9").Value = '  +6.77%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06646'
$ws.Range("D11").Value = '1.892.98'
$ws.Range("E11").Value = '  +4.00%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '17.01'
$ws.Range("E12").Value = '  +2.29%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07226'
$ws.Range("E13").Value = '  +2.50%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6777'
$ws.Range("E14").Value = '  +5.68%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '85.92'
$ws.Range("E15").Value = '  +2.31%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.852'
$ws.Range("E16").Value = '  +3.74%  '
$ws.Range("D17").Value = '30.042.88'
$ws.Range("E17").Value = '  +3.65%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000007979'
$ws.Range("E18").Value = '  +9.51%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9978'
$ws.Range("E19").Value = '  -0.23%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.94'
$ws.Range("E20").Value = '  +6.20%  '
$ws.Range("D21").Value = '2.134.50'
$ws.Range("E21").Value = '  +4.04%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9978'
$ws.Range("E22").Value = '  -0.34%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.781'
$ws.Range("E23").Value = '  +5.39%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.690'
$ws.Range("E24").Value = '  +6.38%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.206'
$ws.Range("E25").Value = '  +4.58%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '147.45'
$ws.Range("E26").Value = '  +2.33%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '131.90'
$ws.Range("E27").Value = '  +3.28%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.85'
$ws.Range("E28").Value = '  +3.32%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.956'
$ws.Range("E29").Value = '  +4.41%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.366'
$ws.Range("E30").Value = '  -2.30%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.256'
$ws.Range("E31").Value = '  +3.50%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08758'
$ws.Range("E32").Value = '  +5.07%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.958'
$ws.Range("E33").Value = '  +5.17%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05095'
$ws.Range("E34").Value = '  +3.14%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.121'
$ws.Range("E35").Value = '  +2.48%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7054'
$ws.Range("E36").Value = '  +5.55%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.665'
$ws.Range("E37").Value = '  -0.94%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.779'
$ws.Range("E38").Value = '  +3.68%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.224'
$ws.Range("E39").Value = '  -2.68%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9440'
$ws.Range("E40").Value = '  -0.01%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.01660'
$ws.Range("E41").Value = '  +4.92%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.078'
$ws.Range("E42").Value = '  -1.52%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9969'
$ws.Range("E43").Value = '  -0.40%  '
# Row 44/45: coins swapped (TheSandbox now row44, Quant now row45) with updated prices
$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4220'
$ws.Range("E44").Value = '  +4.12%  '
$ws.Range("B45").Value = 'Quant'
$ws.Range("C45").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '103.20'
$ws.Range("E45").Value = '  +2.40%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.493'
$ws.Range("E46").Value = '  +5.09%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1261'
$ws.Range("E47").Value = '  +3.67%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.05724'
$ws.Range("E48").Value = '  +3.64%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '32.81'
$ws.Range("E49").Value = '  +3.77%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.237'
$ws.Range("E50").Value = '  +1.89%  '
$ws.Range("E51").Value = '  +4.34%  '
